$p = $ppt.ActivePresentation

# --- Slide 20: merge the "This is ___ pencil" runs into a single run ---
$s20 = $p.Slides.Item(20)
$tf20 = $s20.Shapes.Item(2).TextFrame
$tr20 = $tf20.TextRange

$para1 = $tr20.Paragraphs(1, 1)
$para1.Text = "PLACEHOLDER_RESET_1"
$para1 = $tr20.Paragraphs(1, 1)
$para1.Text = "This is ___ pencil. (you, your, you’re, are)"

# --- Slide 44: "and" -> "but" and "Obnoxious" -> "Ugly" ---
$s44 = $p.Slides.Item(44)
$tf44 = $s44.Shapes.Item(2).TextFrame
$tr44 = $tf44.TextRange

$paraA = $tr44.Paragraphs(1, 1)
$paraA.Text = "PLACEHOLDER_RESET_2"
$paraA = $tr44.Paragraphs(1, 1)
$paraA.Text = "He is not pretty but not ____. (Obnoxious, Approximation, Book, But)"

$paraB = $tr44.Paragraphs(4, 1)
$paraB.Text = "PLACEHOLDER_RESET_3"
$paraB = $tr44.Paragraphs(4, 1)
$paraB.Text = "I don’t like her. She is pretty but ____. (Ugly, Condescending, Approximation, Attractive)"
